# Atualizei dados bibi e add
# - Reordena os nomes das lojas (Bibi Cell Mundi passa a ficar depois de
#   "Bibi Cell Manauara" em vez de logo apos o cabecalho "nome").
# - Atualiza os valores diarios de faturamento das linhas 2 a 6 (novos dados).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Atualiza os nomes das lojas na coluna A (linhas 2 a 4 trocam de lugar) ---
$ws.Range("A2").Value = "Bibi Cell Vieiralves"
$ws.Range("A3").Value = "Bibi Cell Manauara"
$ws.Range("A4").Value = "Bibi Cell Mundi"
# Linhas 5 e 6 permanecem iguais (Bibi Cell Ponta Negra / total)

# --- Novos valores de faturamento diario (colunas B..AG) por linha ---
# Cada array cobre as colunas B..AG (31 dias + coluna "total" em AG).
$row2 = @(8802,7274,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,16076)
$row3 = @(2469.75,5177,3030,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,10676.75)
$row4 = @(8258,2278,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,10536)
$row5 = @(2221.85,6006.7,1990,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,10218.55)
$row6 = @(21751.6,20735.7,5020,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,47507.3)

$rowNumbers = @(2, 3, 4, 5, 6)
$rowValues = @($row2, $row3, $row4, $row5, $row6)

for ($i = 0; $i -lt $rowNumbers.Length; $i++) {
    $r = $rowNumbers[$i]
    $values = $rowValues[$i]
    $col = 2   # coluna B
    foreach ($v in $values) {
        $ws.Cells.Item($r, $col).Value = $v
        $col = $col + 1
    }
}
